# Weekly data refresh: shift each record down one row, with the
# oldest (row 5) values wrapping around to become the newest (row 2).
# Only the observation columns change (D, J-Q); the descriptive
# columns (A,B,C,E-I,R) are identical across all rows already.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the "before" values for rows 2-5, columns D and J:Q
$rows = 2..5
$before = @{}
foreach ($r in $rows) {
    $before[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        J = $ws.Cells.Item($r, 10).Value2
        K = $ws.Cells.Item($r, 11).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        N = $ws.Cells.Item($r, 14).Value2
        O = $ws.Cells.Item($r, 15).Value2
        P = $ws.Cells.Item($r, 16).Value2
        Q = $ws.Cells.Item($r, 17).Value2
    }
}

# New row r gets the old values of row (r-1), wrapping row 2 <- row 5
$srcFor = @{ 2 = 5; 3 = 2; 4 = 3; 5 = 4 }

foreach ($r in $rows) {
    $src = $before[$srcFor[$r]]
    $ws.Cells.Item($r, 4).Value2  = $src.D
    $ws.Cells.Item($r, 10).Value2 = $src.J
    $ws.Cells.Item($r, 11).Value2 = $src.K
    $ws.Cells.Item($r, 12).Value2 = $src.L
    $ws.Cells.Item($r, 13).Value2 = $src.M
    $ws.Cells.Item($r, 14).Value2 = $src.N
    $ws.Cells.Item($r, 15).Value2 = $src.O
    $ws.Cells.Item($r, 16).Value2 = $src.P
    $ws.Cells.Item($r, 17).Value2 = $src.Q
}
